# Add new job posting row JD_005 (Test Engineer / Test, min 2 max 3 years)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "JD_005"
$ws.Range("B6").Value = "Test Engineer"
$ws.Range("C6").Value = "Test"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 3
